# Update the "Goals" bullet on the Requirements/Scope slide (slide 3) so
# that it reads "... Session-Reflector - stateless mode" instead of
# "... Session-Reflector for stateless mode".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$targetOld = "Avoid provisioning and maintaining test sessions on Session-Reflector for stateless mode"
$targetNew = "Avoid provisioning and maintaining test sessions on Session-Reflector - stateless mode"

$paragraphCount = $textRange.Paragraphs().Count
for ($i = 1; $i -le $paragraphCount; $i++) {
    $para = $textRange.Paragraphs($i)
    # Paragraph.Text includes a trailing paragraph-mark character (CR),
    # so trim it before comparing against the plain target string.
    if ($para.Text.TrimEnd("`r") -eq $targetOld) {
        # Select the whole paragraph's text as one characters range so the
        # run is replaced in-place (keeping a single run with its existing
        # formatting) rather than being split into multiple runs.
        $whole = $para.Characters(1, $para.Text.Length)
        $whole.Text = $targetNew
    }
}
